## Fruta / hortaliza, semanal
##
## Weekly refresh of the price series: a new observation is inserted as
## row 5 (right after the 3 most-recent weeks already at the top of the
## table) and every older observation shifts down one row - so the row
## that used to be the last one (42) ends up preserved as the new last
## row (43) of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 5..42 down one slot, opening up a blank row 5 (inherits the
# Fecha column's date formatting from the surrounding rows automatically,
# same as Excel's own Insert Row).
$ws.Rows.Item(5).Insert()

# Fill the new weekly observation into row 5.
$ws.Cells.Item(5, 1).Value  = 2
$ws.Cells.Item(5, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(5, 3).Value  = "Coquimbo"
$ws.Cells.Item(5, 4).Value  = 44490
$ws.Cells.Item(5, 5).Value  = 4
$ws.Cells.Item(5, 6).Value  = 100112026
$ws.Cells.Item(5, 7).Value  = "Haba"
$ws.Cells.Item(5, 8).Value  = "Sin especificar"
$ws.Cells.Item(5, 9).Value  = "Primera"
$ws.Cells.Item(5, 10).Value = 400
$ws.Cells.Item(5, 11).Value = 5000
$ws.Cells.Item(5, 12).Value = 6000
$ws.Cells.Item(5, 13).Value = 5500
$ws.Cells.Item(5, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(5, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(5, 16).Value = 220
$ws.Cells.Item(5, 17).Value = 25
$ws.Cells.Item(5, 18).Value = "Hortaliza"
